# "try to fix render markdown"
#
# The "along" column (E) uses the label "<b>All</b>" (HTML bold, meant to be
# rendered as markdown) for the first row of every 4-row block. Replace that
# label with the markdown-style "*All*" so it renders correctly wherever this
# sheet's HTML-ish bold tags were not being interpreted.
#
# Note: column A's "<b>All</b>" (the region/country label) is left untouched -
# only the vote-category label in column E is affected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 2..45
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 5)
    $current = $cell.Value()
    if ($current -eq "<b>All</b>") {
        $cell.Value = "*All*"
    }
}
